# ---------------------------------------------------------------------------
# "Added analysis for close types"
#
# The original sheet had a single 2-row table (Females/Males) of 15 trait
# percentages with a line chart. This edit:
#   1. Inserts a header row above the table (trait names) and adds it as a
#      chart category axis.
#   2. Adds a small stats scratch table (topic/effect -> AVERAGE/STDEV.P).
#   3. Adds a second, re-ordered ("sorted to communal/agentic") copy of the
#      table further down the sheet, with its own line chart.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Shift the existing 2-row table down by inserting a new row 1, then
#    populate it with the trait-name header used as the chart category axis.
#    (I1 is intentionally left for later so the shared-string order matches.)
# ---------------------------------------------------------------------------
$ws.Rows("1").Insert()

$ws.Range("B1").Value2 = "Determined"
$ws.Range("C1").Value2 = "Collaborator"
$ws.Range("D1").Value2 = "Family"
$ws.Range("E1").Value2 = "Optimistic"
$ws.Range("F1").Value2 = "Conqueror"
$ws.Range("G1").Value2 = "Peacful zionist"
$ws.Range("H1").Value2 = "Nurturing"

$ws.Range("J1").Value2 = "Educator"
$ws.Range("K1").Value2 = "Army Commander"
$ws.Range("L1").Value2 = "Visionary"
$ws.Range("M1").Value2 = "N/A"
$ws.Range("N1").Value2 = "Achiever"
$ws.Range("O1").Value2 = "Fiction Character"
$ws.Range("P1").Value2 = "Problem Sover"

# ---------------------------------------------------------------------------
# 2. Small stats scratch table: topic / effect, with AVERAGE & STDEV.P.
# ---------------------------------------------------------------------------
$ws.Range("C6").Value2 = "topic"
$ws.Range("D6").Value2 = "effect"

$ws.Range("C7").Value2 = 3
$ws.Range("D7").Value2 = 9.95
$ws.Range("C8").Value2 = 4
$ws.Range("D8").Value2 = 32.48
$ws.Range("C9").Value2 = 6
$ws.Range("D9").Value2 = 19.86
$ws.Range("C10").Value2 = 8
$ws.Range("D10").Value2 = 55.47
$ws.Range("C11").Value2 = 10
$ws.Range("D11").Value2 = 14.28

$ws.Range("F7").Formula = "=AVERAGE(D7:D11)"
$ws.Range("F8").Formula = "=STDEV.P(D7:D11)"

# ---------------------------------------------------------------------------
# 3. Second table label + header row (re-ordered "communal/agentic" sort).
# ---------------------------------------------------------------------------
$ws.Range("A14").Value2 = "Sorted to communal/agentic"

# Now fill I1 with "Patronistic" (this is the last *new* unique string to
# enter the shared-string table from row 1).
$ws.Range("I1").Value2 = "Patronistic"

$ws.Range("B15").Value2 = "Nurturing"
$ws.Range("C15").Value2 = "Collaborator"
$ws.Range("D15").Value2 = "Family"
$ws.Range("E15").Value2 = "Optimistic"
$ws.Range("F15").Value2 = "Fiction Character"
$ws.Range("G15").Value2 = "Educator"
$ws.Range("H15").Value2 = "Visionary"
$ws.Range("I15").Value2 = "Determined"
$ws.Range("J15").Value2 = "Conqueror"
$ws.Range("K15").Value2 = "Peacful zionist"
$ws.Range("L15").Value2 = "Patronistic"
$ws.Range("M15").Value2 = "Army Commander"
$ws.Range("N15").Value2 = "Achiever"
$ws.Range("O15").Value2 = "Problem Sover"

# Header rows get an explicit black font (matches the new style applied in
# the workbook: font 18 / cellXf 1).
$ws.Range("B1:P1").Font.Color = 0
$ws.Range("B15:O15").Font.Color = 0

# ---------------------------------------------------------------------------
# 4. Re-ordered data rows (16/17): same percentages already present in rows
#    2/3, just reshuffled into the column order defined by row 15 (copied
#    value-for-value off the original cells so precision is exact and the
#    N/A column is simply dropped).
# ---------------------------------------------------------------------------
$ws.Range("A16").Value2 = "Females"
$ws.Range("B16").Value2 = $ws.Range("H2").Value2   # Nurturing
$ws.Range("C16").Value2 = $ws.Range("C2").Value2   # Collaborator
$ws.Range("D16").Value2 = $ws.Range("D2").Value2   # Family
$ws.Range("E16").Value2 = $ws.Range("E2").Value2   # Optimistic
$ws.Range("F16").Value2 = $ws.Range("O2").Value2   # Fiction Character
$ws.Range("G16").Value2 = $ws.Range("J2").Value2   # Educator
$ws.Range("H16").Value2 = $ws.Range("L2").Value2   # Visionary
$ws.Range("I16").Value2 = $ws.Range("B2").Value2   # Determined
$ws.Range("J16").Value2 = $ws.Range("F2").Value2   # Conqueror
$ws.Range("K16").Value2 = $ws.Range("G2").Value2   # Peacful zionist
$ws.Range("L16").Value2 = $ws.Range("I2").Value2   # Patronistic
$ws.Range("M16").Value2 = $ws.Range("K2").Value2   # Army Commander
$ws.Range("N16").Value2 = $ws.Range("N2").Value2   # Achiever
$ws.Range("O16").Value2 = $ws.Range("P2").Value2   # Problem Sover

$ws.Range("A17").Value2 = "Males"
$ws.Range("B17").Value2 = $ws.Range("H3").Value2   # Nurturing
$ws.Range("C17").Value2 = $ws.Range("C3").Value2   # Collaborator
$ws.Range("D17").Value2 = $ws.Range("D3").Value2   # Family
$ws.Range("E17").Value2 = $ws.Range("E3").Value2   # Optimistic
$ws.Range("F17").Value2 = $ws.Range("O3").Value2   # Fiction Character
$ws.Range("G17").Value2 = $ws.Range("J3").Value2   # Educator
$ws.Range("H17").Value2 = $ws.Range("L3").Value2   # Visionary
$ws.Range("I17").Value2 = $ws.Range("B3").Value2   # Determined
$ws.Range("J17").Value2 = $ws.Range("F3").Value2   # Conqueror
$ws.Range("K17").Value2 = $ws.Range("G3").Value2   # Peacful zionist
$ws.Range("L17").Value2 = $ws.Range("I3").Value2   # Patronistic
$ws.Range("M17").Value2 = $ws.Range("K3").Value2   # Army Commander
$ws.Range("N17").Value2 = $ws.Range("N3").Value2   # Achiever
$ws.Range("O17").Value2 = $ws.Range("P3").Value2   # Problem Sover

# ---------------------------------------------------------------------------
# 5. Fix up the existing chart (now references rows 2/3 instead of 1/2, and
#    gets trait names as its category axis).
# ---------------------------------------------------------------------------
$co1 = $ws.ChartObjects().Item(1)
$chart1 = $co1.Chart

$s1 = $chart1.SeriesCollection().Item(1)
$s1.Formula = "=SERIES(males_and_females_on_females_me!`$A`$2,males_and_females_on_females_me!`$B`$1:`$P`$1,males_and_females_on_females_me!`$B`$2:`$P`$2,1)"

$s2 = $chart1.SeriesCollection().Item(2)
$s2.Formula = "=SERIES(males_and_females_on_females_me!`$A`$3,males_and_females_on_females_me!`$B`$1:`$P`$1,males_and_females_on_females_me!`$B`$3:`$P`$3,2)"

# Move it further down the sheet to make room for the new chart above it.
$co1.Top = 255
$co1.Left = 935
$co1.Width = 409
$co1.Height = 225

# ---------------------------------------------------------------------------
# 6. Add the second chart for the re-ordered "communal/agentic" table.
# ---------------------------------------------------------------------------
$co2 = $ws.ChartObjects().Add(935, 30, 409, 210)
$co2.Name = "Chart 2"
$chart2 = $co2.Chart
$chart2.ChartType = 65

$ns1 = $chart2.SeriesCollection().NewSeries()
$ns1.Name = "=males_and_females_on_females_me!`$A`$16"
$ns1.XValues = "=males_and_females_on_females_me!`$B`$15:`$O`$15"
$ns1.Values = "=males_and_females_on_females_me!`$B`$16:`$O`$16"

$ns2 = $chart2.SeriesCollection().NewSeries()
$ns2.Name = "=males_and_females_on_females_me!`$A`$17"
$ns2.XValues = "=males_and_females_on_females_me!`$B`$15:`$O`$15"
$ns2.Values = "=males_and_females_on_females_me!`$B`$17:`$O`$17"

$chart2.HasLegend = $true
$chart2.Legend.Position = -4107

# ---------------------------------------------------------------------------
# 7. View tweaks (zoom + active cell) matching the end-state sheet view.
# ---------------------------------------------------------------------------
$ws.Range("F9").Select()
$excel.ActiveWindow.Zoom = 80
